# ET kalender 18-19 06102018.xlsx -- "Kalender bijgewerkt" update
#
# Summary of changes:
#  - Reitse (FotoFinish) renamed to Floris Klunder; a second FotoFinish
#    person, Karuud Pots, is added in a new column G for the NK Junioren
#    rows.
#  - Several manager/timer assignments in the calendar are corrected /
#    swapped / cleared for a handful of rows.
#  - The summary COUNTIF table (rows 38-47) recalculates automatically
#    from those changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Rename the "Reitse" FotoFinish person to "Floris Klunder" -------------
# (F23/F24/F25 already hold this shared string; changing the text updates
# all three cells at once, matching the sharedStrings.xml edit.)
$ws.Range("F23").Value = "Floris Klunder"
$ws.Range("F24").Value = "Floris Klunder"
$ws.Range("F25").Value = "Floris Klunder"

# --- Add the new FotoFinish column (G) --------------------------------------
$ws.Range("G1").Value = $ws.Range("F1").Value()

$ws.Range("F23").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("G23").Value = "Karuud Pots"

$ws.Range("F24").Copy()
$ws.Range("G24").PasteSpecial(-4122)
$ws.Range("G24").Value = "Karuud Pots"

$ws.Range("F25").Copy()
$ws.Range("G25").PasteSpecial(-4122)
$ws.Range("G25").Value = "Karuud Pots"

$ws.Range("F1:G1").EntireColumn.AutoFit()

# --- Manager / timer corrections in the calendar rows -----------------------
# Row 13 (Grunobokaal): Timer2 Jan Reijne -> Reitse Eskens
$ws.Range("E13").Value = "Reitse Eskens"

# Row 20 (Klassewedstrijd): Timer Evelien Eskens removed
$ws.Range("D20").ClearContents()
$ws.Range("D20").Interior.Color = 65535

# Row 23 (NK Junioren): Manager Richard Verschure -> Reitse Eskens
$ws.Range("C23").Value = "Reitse Eskens"

# Row 24 (NK Junioren): Manager Peter Nomden -> Reitse Eskens; Timer2 Jan Reijne removed
$ws.Range("C24").Value = "Reitse Eskens"
$ws.Range("E24").ClearContents()
$ws.Range("E24").Interior.Color = 65535

# Row 25 (NK Junioren): Manager Rudy Kok -> Reitse Eskens
$ws.Range("C25").Value = "Reitse Eskens"

# Row 26 / 27 (Klassewedstrijd / Baankampioenschap afstanden): Timer swapped
$ws.Range("D26").Value = "Gea Kunst"
$ws.Range("D27").Value = "Gerrit Visser"

# Row 28 (Baankampioenschap afstanden dag 2): Timer Jan Reijne removed
$ws.Range("D28").ClearContents()
$ws.Range("D28").Interior.Color = 65535

# Row 29 (Baankampioenschappen jeugd): Manager Gerrit Visser -> Gea Kunst
$ws.Range("C29").Value = "Gea Kunst"

# Row 30 (Klassewedstrijd): Timer Jan Reijne removed
$ws.Range("D30").ClearContents()
$ws.Range("D30").Interior.Color = 65535

# --- Move the saved selection to D34, with no frozen scroll position --------
$ws.Range("D34").Select()
